# Refresh the crypto symbol list (prices, 1h volume %, and the "Hora"
# hour marker) to match the GitHub Actions scrape taken at 03:08 UTC.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These columns are stored as plain text in the workbook (prices such as
# "295.60" and percentages such as "1.89%" are not real numbers). Mark the
# cells we are about to touch as Text first so Excel does not silently
# reinterpret them as numbers/percentages when we assign the new values.
$ws.Range("D2:D11").NumberFormat = "@"
$ws.Range("D13:D18").NumberFormat = "@"
$ws.Range("D21:D26").NumberFormat = "@"
$ws.Range("D39:D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"
$ws.Range("E39:E41").NumberFormat = "@"
$ws.Range("E44:E47").NumberFormat = "@"
$ws.Range("E49:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# New values scraped for each coin row (Price / Volume(1h) / Hora).
$ws.Range("D2").Value = "295.60"
$ws.Range("E2").Value = "1.89%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "31.20"
$ws.Range("E3").Value = "1.13%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "4.946"
$ws.Range("E4").Value = "0.07%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.07447"
$ws.Range("E5").Value = "4.31%"
$ws.Range("G5").Value = "3"
$ws.Range("D6").Value = "2.133"
$ws.Range("E6").Value = "18.18%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "7.760"
$ws.Range("E7").Value = "1.01%"
$ws.Range("G7").Value = "3"
$ws.Range("D8").Value = "3.747"
$ws.Range("E8").Value = "0.37%"
$ws.Range("G8").Value = "3"
$ws.Range("D9").Value = "0.9152"
$ws.Range("E9").Value = "2.10%"
$ws.Range("G9").Value = "3"
$ws.Range("D10").Value = "0.08673"
$ws.Range("E10").Value = "14.78%"
$ws.Range("G10").Value = "3"
$ws.Range("D11").Value = "0.1699"
$ws.Range("E11").Value = "3.27%"
$ws.Range("G11").Value = "3"
$ws.Range("E12").Value = "2.81%"
$ws.Range("G12").Value = "3"
$ws.Range("D13").Value = "0.03151"
$ws.Range("E13").Value = "3.33%"
$ws.Range("G13").Value = "3"
$ws.Range("D14").Value = "0.1007"
$ws.Range("E14").Value = "0.55%"
$ws.Range("G14").Value = "3"
$ws.Range("D15").Value = "0.001509"
$ws.Range("E15").Value = "0.60%"
$ws.Range("G15").Value = "3"
$ws.Range("D16").Value = "0.005812"
$ws.Range("E16").Value = "-0.30%"
$ws.Range("G16").Value = "3"
$ws.Range("D17").Value = "3.512"
$ws.Range("E17").Value = "1.06%"
$ws.Range("G17").Value = "3"
$ws.Range("D18").Value = "2.078"
$ws.Range("E18").Value = "-0.08%"
$ws.Range("G18").Value = "3"
$ws.Range("E19").Value = "1.61%"
$ws.Range("G19").Value = "3"
$ws.Range("E20").Value = "-0.16%"
$ws.Range("G20").Value = "3"
$ws.Range("D21").Value = "3.973"
$ws.Range("E21").Value = "-1.53%"
$ws.Range("G21").Value = "3"
$ws.Range("D22").Value = "0.2102"
$ws.Range("E22").Value = "5.22%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04557"
$ws.Range("E23").Value = "1.15%"
$ws.Range("G23").Value = "3"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").Value = "0.19%"
$ws.Range("G24").Value = "3"
$ws.Range("D25").Value = "0.004623"
$ws.Range("E25").Value = "15.75%"
$ws.Range("G25").Value = "3"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "4.14%"
$ws.Range("G26").Value = "3"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "0.00%"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.04480"
$ws.Range("E40").Value = "2.85%"
$ws.Range("G40").Value = "3"
$ws.Range("D41").Value = "0.007278"
$ws.Range("E41").Value = "-1.30%"
$ws.Range("G41").Value = "3"
$ws.Range("D42").Value = "0.008992"
$ws.Range("G42").Value = "3"
$ws.Range("D43").Value = "0.1331"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.001971"
$ws.Range("E44").Value = "-1.38%"
$ws.Range("G44").Value = "3"
$ws.Range("D45").Value = "0.009132"
$ws.Range("E45").Value = "-2.93%"
$ws.Range("G45").Value = "3"
$ws.Range("D46").Value = "0.00006103"
$ws.Range("E46").Value = "1.33%"
$ws.Range("G46").Value = "3"
$ws.Range("E47").Value = "0.12%"
$ws.Range("G47").Value = "3"
$ws.Range("G48").Value = "3"
$ws.Range("D49").Value = "0.002002"
$ws.Range("E49").Value = "-33.22%"
$ws.Range("G49").Value = "3"
$ws.Range("E50").Value = "0.12%"
$ws.Range("G50").Value = "3"
$ws.Range("E51").Value = "0.12%"
$ws.Range("G51").Value = "3"

Write-Host "Updated 116 cells."
